$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("733:736").Insert()

# Row 733
$ws.Cells.Item(733, 1).Value = 3
$ws.Cells.Item(733, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(733, 3).Value = 'Coquimbo'
$ws.Cells.Item(733, 4).Value = 44585
$ws.Cells.Item(733, 5).Value = 5
$ws.Cells.Item(733, 6).Value = 'Fruta'
$ws.Cells.Item(733, 7).Value = 100103
$ws.Cells.Item(733, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(733, 9).Value = 100103004
$ws.Cells.Item(733, 10).Value = 'Durazno'
$ws.Cells.Item(733, 11).Value = 'Carson'
$ws.Cells.Item(733, 12).Value = 'Especial'
$ws.Cells.Item(733, 13).Value = 80
$ws.Cells.Item(733, 14).Value = 16000
$ws.Cells.Item(733, 15).Value = 16000
$ws.Cells.Item(733, 16).Value = 16000
$ws.Cells.Item(733, 17).Value = '$/caja 15 kilos empedrada'
$ws.Cells.Item(733, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(733, 19).Value = 1067
$ws.Cells.Item(733, 20).Value = 15

# Row 734
$ws.Cells.Item(734, 1).Value = 3
$ws.Cells.Item(734, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(734, 3).Value = 'Coquimbo'
$ws.Cells.Item(734, 4).Value = 44585
$ws.Cells.Item(734, 5).Value = 5
$ws.Cells.Item(734, 6).Value = 'Fruta'
$ws.Cells.Item(734, 7).Value = 100103
$ws.Cells.Item(734, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(734, 9).Value = 100103004
$ws.Cells.Item(734, 10).Value = 'Durazno'
$ws.Cells.Item(734, 11).Value = 'Carson'
$ws.Cells.Item(734, 12).Value = 'Extra (doble especial)'
$ws.Cells.Item(734, 13).Value = 50
$ws.Cells.Item(734, 14).Value = 18000
$ws.Cells.Item(734, 15).Value = 18000
$ws.Cells.Item(734, 16).Value = 18000
$ws.Cells.Item(734, 17).Value = '$/caja 15 kilos empedrada'
$ws.Cells.Item(734, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(734, 19).Value = 1200
$ws.Cells.Item(734, 20).Value = 15

# Row 735
$ws.Cells.Item(735, 1).Value = 3
$ws.Cells.Item(735, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(735, 3).Value = 'Coquimbo'
$ws.Cells.Item(735, 4).Value = 44585
$ws.Cells.Item(735, 5).Value = 5
$ws.Cells.Item(735, 6).Value = 'Fruta'
$ws.Cells.Item(735, 7).Value = 100103
$ws.Cells.Item(735, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(735, 9).Value = 100103004
$ws.Cells.Item(735, 10).Value = 'Durazno'
$ws.Cells.Item(735, 11).Value = 'Carson'
$ws.Cells.Item(735, 12).Value = 'Primera'
$ws.Cells.Item(735, 13).Value = 150
$ws.Cells.Item(735, 14).Value = 14000
$ws.Cells.Item(735, 15).Value = 14500
$ws.Cells.Item(735, 16).Value = 14267
$ws.Cells.Item(735, 17).Value = '$/caja 15 kilos empedrada'
$ws.Cells.Item(735, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(735, 19).Value = 951
$ws.Cells.Item(735, 20).Value = 15

# Row 736
$ws.Cells.Item(736, 1).Value = 3
$ws.Cells.Item(736, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(736, 3).Value = 'Coquimbo'
$ws.Cells.Item(736, 4).Value = 44560
$ws.Cells.Item(736, 5).Value = 5
$ws.Cells.Item(736, 6).Value = 'Fruta'
$ws.Cells.Item(736, 7).Value = 100103
$ws.Cells.Item(736, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(736, 9).Value = 100103004
$ws.Cells.Item(736, 10).Value = 'Durazno'
$ws.Cells.Item(736, 11).Value = 'Carson'
$ws.Cells.Item(736, 12).Value = 'Segunda'
$ws.Cells.Item(736, 13).Value = 90
$ws.Cells.Item(736, 14).Value = 12000
$ws.Cells.Item(736, 15).Value = 12000
$ws.Cells.Item(736, 16).Value = 12000
$ws.Cells.Item(736, 17).Value = '$/caja 15 kilos empedrada'
$ws.Cells.Item(736, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(736, 19).Value = 800
$ws.Cells.Item(736, 20).Value = 15
